$d = $word.ActiveDocument

# Locate the target paragraph: "PROBLEM: pathfinding costs for traversing south east seem to be broken :D"
$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*pathfinding costs for traversing south east seem to be broken*") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)

# 1) Insert a brand new paragraph BEFORE the target paragraph, inheriting the
#    same paragraph formatting (ListParagraph style, numPr ilvl 0 / numId 4, rPr lang en-GB),
#    and give it the new "PROBLEM: I ran into a bug..." text.
$target.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($targetIndex)
$newPara.Range.Text = "PROBLEM: I ran into a bug where if there was lag, the player would continually walk in a certain direction. I resolved this by updating the direction the player should travel in before moving the player in that direction as opposed to only calculating the direction once."

# 2) Re-locate the pathfinding paragraph (it has shifted down by one) and append a
#    new run containing the solution text after the existing "PROBLEM: pathfinding..." run.
$targetIndex2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*pathfinding costs for traversing south east seem to be broken*") {
        $targetIndex2 = $i
        break
    }
}

$target2 = $d.Paragraphs.Item($targetIndex2)
$r = $target2.Range
# Exclude the paragraph mark from the range (End - 1) so the insertion happens
# at the end of the existing text, before the pilcrow.
$insertPoint = $d.Range($r.End - 1, $r.End - 1)
$insertPoint.InsertAfter(" – Solution: When calculating the distance cost, I was incorrectly calculating the difference of the y because I was using the x variable for node a…..")
# Toggling formatting forces the newly inserted text to remain its own run
# rather than being silently re-merged into the preceding identical run.
$insertPoint.Bold = 1
$insertPoint.Bold = 0
